# Append: 2025-09-09 06:26 JST
# Refresh the "ランサーズ" listing sheet: drop the old rows 9-18, replace the
# remaining rows 2-8 with the newly scraped listings, and resize two columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Remove the hyperlinks first (they would otherwise keep pointing at
#        cells in rows that are about to be deleted / retargeted). ---
$ws.Hyperlinks.Delete()

# --- 2. Drop rows 9-18 entirely (dimension shrinks to A1:H8). ---
$ws.Range("A9:A18").EntireRow.Delete()

# --- 3. Resize columns B and H. ColumnWidth's stored character width is
#        offset by 5/6 from the saved OOXML <col width> value, so subtract
#        that padding to land on the exact target width. ---
$ws.Columns.Item(2).ColumnWidth = 47 - (5/6)
$ws.Columns.Item(8).ColumnWidth = 14 - (5/6)

# --- 4. Rewrite the data rows with the freshly fetched listings. ---
$timestamp = "2025-09-09 06:26:28"

$ws.Range("A2").Value = $timestamp
$ws.Range("B2").Value = "【AI活用】社内業務改善・RAG検索システム構築の依頼"
$ws.Range("D2").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5390045"
$ws.Range("G2").Value = 353
$ws.Range("H2").Value = "🔥AI,Ai ◇業務改善"

$ws.Range("A3").Value = $timestamp
$ws.Range("B3").Value = "【AIクローン制作】肉声・映像応答のプロデュース依頼"
$ws.Range("D3").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5389313"
$ws.Range("G3").Value = 303
$ws.Range("H3").Value = "🔥AI,Ai"

$ws.Range("A4").Value = $timestamp
$ws.Range("B4").Value = "【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪"
$ws.Range("D4").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5217096"
$ws.Range("G4").Value = 243
$ws.Range("H4").Value = "🔥API ◆ツール"

$ws.Range("A5").Value = $timestamp
$ws.Range("B5").Value = "【急募】共同開発エンジニアメンバーを募集します!"
$ws.Range("D5").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5390021"
$ws.Range("G5").Value = 68
$ws.Range("H5").Value = "◆開発"

$ws.Range("A6").Value = $timestamp
$ws.Range("B6").Value = "【仮想通貨】自動売買EAの開発依頼"
$ws.Range("D6").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5389714"
$ws.Range("G6").Value = 63
$ws.Range("H6").Value = "◆開発"

$ws.Range("A7").Value = $timestamp
$ws.Range("B7").Value = "お歳暮受注発注集計スプレッドシート開発依頼"
$ws.Range("D7").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5389953"
$ws.Range("G7").Value = 60
$ws.Range("H7").Value = "◆開発"

$ws.Range("A8").Value = $timestamp
$ws.Range("B8").Value = "限定公開 PR 限定公開の仕事"
$ws.Range("D8").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5385681"
$ws.Range("G8").Value = 25
$ws.Range("H8").ClearContents()

# --- 5. Re-create the hyperlinks on F2:F8 against the refreshed URLs. ---
$ws.Hyperlinks.Add($ws.Range("F2"), $ws.Range("F2").Text)
$ws.Hyperlinks.Add($ws.Range("F3"), $ws.Range("F3").Text)
$ws.Hyperlinks.Add($ws.Range("F4"), $ws.Range("F4").Text)
$ws.Hyperlinks.Add($ws.Range("F5"), $ws.Range("F5").Text)
$ws.Hyperlinks.Add($ws.Range("F6"), $ws.Range("F6").Text)
$ws.Hyperlinks.Add($ws.Range("F7"), $ws.Range("F7").Text)
$ws.Hyperlinks.Add($ws.Range("F8"), $ws.Range("F8").Text)
